$d = $word.ActiveDocument

# Replace the placeholder NIP "-" with the actual NIP number in both
# occurrences of "NIP. -" in the document.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("NIP. -", $true, $false, $false, $false, $false, $true, 1, $false, "NIP. 19851027 201706 1 001", 2)
